# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update first "Bad Drivers" row (row 3) ---
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 98.7

# --- Remove the second "Bad Drivers" data row (old row 4, "23.90.0.2") ---
# This shifts the "Totals:" row (old row 5) up to row 4, and shifts every
# row below it up by one as well.
$ws.Rows("4").Delete()

# --- Update the "Totals:" row (now row 4) ---
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 13

# --- Update "Good Drivers" table (now starting row 10) ---
# Row 12: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1
$ws.Range("B12").Value = 11140

# Row 13: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3
$ws.Range("B13").Value = 14487

Write-Host "edits applied"
